# Apply weekly Fruta / hortaliza price updates to the "Pera" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44355
$ws.Range("N3").Value = 17000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 17500
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("R3").Value = 'Región Metropolitana'
$ws.Range("S3").Value = 972
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("D4").Value = 44355
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17500
$ws.Range("Q4").Value = '$/caja 18 kilos granel'
$ws.Range("R4").Value = 'Región Metropolitana'
$ws.Range("S4").Value = 972
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("D5").Value = 44474
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 19000
$ws.Range("P5").Value = 18500
$ws.Range("Q5").Value = '$/caja 18 kilos empedrada'
$ws.Range("S5").Value = 1028

# Row 6
$ws.Range("K6").Value = 'Winter Nelis'
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 17000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 17500
$ws.Range("Q6").Value = '$/bandeja 18 kilos granel'
$ws.Range("S6").Value = 972

# Row 7
$ws.Range("D7").Value = 44280
$ws.Range("K7").Value = 'Packham''s Triumph'
$ws.Range("M7").Value = 350
$ws.Range("N7").Value = 24000
$ws.Range("O7").Value = 25000
$ws.Range("P7").Value = 24500
$ws.Range("Q7").Value = '$/caja 18 kilos granel'
$ws.Range("S7").Value = 1361

# Row 8
$ws.Range("D8").Value = 44280
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 24000
$ws.Range("O8").Value = 25000
$ws.Range("P8").Value = 24500
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("S8").Value = 1361

# Row 9
$ws.Range("D9").Value = 44329
$ws.Range("M9").Value = 340
$ws.Range("N9").Value = 21000
$ws.Range("O9").Value = 22000
$ws.Range("P9").Value = 21500
$ws.Range("Q9").Value = '$/bandeja 18 kilos granel'
$ws.Range("S9").Value = 1194

# Row 10
$ws.Range("D10").Value = 44313
$ws.Range("K10").Value = 'Winter Nelis'
$ws.Range("L10").Value = 'Tercera'
$ws.Range("M10").Value = 250
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15500
$ws.Range("S10").Value = 861

# Row 11
$ws.Range("D11").Value = 44341
$ws.Range("K11").Value = 'Packham''s Triumph'
$ws.Range("M11").Value = 300
$ws.Range("Q11").Value = '$/caja 18 kilos granel'
$ws.Range("R11").Value = 'Región Metropolitana'

# Row 12
$ws.Range("D12").Value = 44678

# Row 13
$ws.Range("D13").Value = 44642
$ws.Range("M13").Value = 270
$ws.Range("N13").Value = 19000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 19500
$ws.Range("S13").Value = 1083

# Row 14
$ws.Range("D14").Value = 44371
$ws.Range("L14").Value = 'Calibre 90'
$ws.Range("M14").Value = 140
$ws.Range("N14").Value = 17000
$ws.Range("O14").Value = 18000
$ws.Range("P14").Value = 17429
$ws.Range("Q14").Value = '$/caja 18 kilos embalada'
$ws.Range("S14").Value = 968

# Row 15
$ws.Range("D15").Value = 44371
$ws.Range("K15").Value = 'Winter Nelis'
$ws.Range("L15").Value = 'Calibre 80'
$ws.Range("M15").Value = 120
$ws.Range("N15").Value = 17000
$ws.Range("O15").Value = 18000
$ws.Range("P15").Value = 17500
$ws.Range("Q15").Value = '$/caja 18 kilos embalada'
$ws.Range("S15").Value = 972

# Row 16
$ws.Range("D16").Value = 44292
$ws.Range("K16").Value = 'Packham''s Triumph'
$ws.Range("M16").Value = 300
$ws.Range("N16").Value = 22000
$ws.Range("O16").Value = 23000
$ws.Range("P16").Value = 22500
$ws.Range("Q16").Value = '$/caja 18 kilos granel'
$ws.Range("S16").Value = 1250

# Row 17
$ws.Range("K17").Value = 'Winter Nelis'
$ws.Range("M17").Value = 250

# Row 18
$ws.Range("D18").Value = 44323
$ws.Range("K18").Value = 'Packham''s Triumph'
$ws.Range("N18").Value = 15000
$ws.Range("O18").Value = 16000
$ws.Range("P18").Value = 15500
$ws.Range("Q18").Value = '$/bandeja 18 kilos granel'
$ws.Range("S18").Value = 861

# Row 19
$ws.Range("D19").Value = 44525
$ws.Range("N19").Value = 19000
$ws.Range("O19").Value = 20000
$ws.Range("P19").Value = 19500
$ws.Range("R19").Value = 'Región de O''Higgins'
$ws.Range("S19").Value = 1083

# Row 20
$ws.Range("D20").Value = 44398
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 20000
$ws.Range("O20").Value = 21000
$ws.Range("P20").Value = 20500
$ws.Range("Q20").Value = '$/caja 20 kilos empedrada'
$ws.Range("S20").Value = 1025
$ws.Range("T20").Value = 20

# Row 21
$ws.Range("D21").Value = 44398
$ws.Range("M21").Value = 200
$ws.Range("N21").Value = 20000
$ws.Range("O21").Value = 21000
$ws.Range("P21").Value = 20500
$ws.Range("Q21").Value = '$/caja 20 kilos empedrada'
$ws.Range("S21").Value = 1025
$ws.Range("T21").Value = 20

# Row 22
$ws.Range("D22").Value = 44497
$ws.Range("M22").Value = 300
$ws.Range("Q22").Value = '$/bandeja 18 kilos granel'
$ws.Range("R22").Value = 'Región de O''Higgins'

# Row 23
$ws.Range("D23").Value = 44497
$ws.Range("Q23").Value = '$/bandeja 18 kilos granel'
$ws.Range("R23").Value = 'Región de O''Higgins'

# Row 24
$ws.Range("D24").Value = 44421
$ws.Range("L24").Value = 'Segunda'
$ws.Range("M24").Value = 270
$ws.Range("N24").Value = 16000
$ws.Range("O24").Value = 17000
$ws.Range("P24").Value = 16500
$ws.Range("Q24").Value = '$/bandeja 18 kilos granel'
$ws.Range("S24").Value = 917

# Row 25
$ws.Range("D25").Value = 44421
$ws.Range("L25").Value = 'Segunda'
$ws.Range("M25").Value = 250
$ws.Range("N25").Value = 16000
$ws.Range("O25").Value = 17000
$ws.Range("P25").Value = 16500
$ws.Range("Q25").Value = '$/bandeja 18 kilos granel'
$ws.Range("S25").Value = 917

# Row 26 (new row)
$ws.Range("A26").Value = 1
$ws.Range("B26").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C26").Value = 'Arica y Parinacota'
$ws.Range("D26").Value = 44314
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = 'Fruta'
$ws.Range("G26").Value = 100104
$ws.Range("H26").Value = 'Frutos de pepita'
$ws.Range("I26").Value = 100104005
$ws.Range("J26").Value = 'Pera'
$ws.Range("K26").Value = 'Packham''s Triumph'
$ws.Range("L26").Value = 'Segunda'
$ws.Range("M26").Value = 250
$ws.Range("N26").Value = 17000
$ws.Range("O26").Value = 18000
$ws.Range("P26").Value = 17500
$ws.Range("Q26").Value = '$/caja 18 kilos granel'
$ws.Range("R26").Value = 'Región de O''Higgins'
$ws.Range("S26").Value = 972
$ws.Range("T26").Value = 18
$ws.Range("D26").NumberFormat = "YYYY-MM-DD HH:MM:SS"

